$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected data (bug fix in lowcode.stats row ordering/values).
# Columns: A,B,C,D,E,F
$data = @{
    2  = @(901, 16, 15, 45, 60, 60)
    3  = @(1001, 18, 30, 75, 60, 72)
    4  = @(301, 6, 45, 30, 60, 45)
    5  = @(501, 9, 52, 30, 75, 45)
    6  = @(601, 9, 60, 67, 60, 42)
    7  = @(201, 9, 30, 15, 45, 30)
    8  = @(902, 1, 0, 0, 0, 0)
    9  = @(401, 9, 48, 67, 75, 45)
    10 = @(1202, 2, 10, 10, 10, 10)
    11 = @(1201, 2, 10, 10, 10, 10)
    12 = @(1203, 3, 15, 15, 15, 15)
    13 = @(101, 9, 30, 15, 60, 15)
    14 = @(701, 3, 90, 45, 97, 15)
    15 = @(801, 3, 67, 65, 52, 45)
    16 = @(1, 0, 2, 2, 2, 2)
    17 = @(802, 0, 4, 5, 4, 0)
    20 = @(502, 0, 4, 0, 0, 0)
    21 = @(1101, 0, 15, 30, 30, 0)
}

$cols = @("A", "B", "C", "D", "E", "F")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
